$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '42.581.47'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.304.19'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'316.52"
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").Value = "'103.38"
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = "'0.605"
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").Value = "'39.64"
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = "'8.49"
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("E13").Value = '  +1.58%  '
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").Value = "'15.36"
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '2.653.79'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '2.299.06'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '42.568.66'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = "'7.57"
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").Value = "'14.12"
$ws.Range("E20").Value = '  +27.30%  '
$ws.Range("D21").Value = "'0.0000106"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = "'74.04"
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").Value = "'3.55"
$ws.Range("E23").Value = '  -3.07%  '
$ws.Range("D24").Value = "'266.98"
$ws.Range("E24").Value = '  -4.85%  '
$ws.Range("E25").Value = '  -2.02%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = '  -1.49%  '
$ws.Range("D29").Value = "'22.71"
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").Value = "'6.73"
$ws.Range("E30").Value = '  +12.95%  '
$ws.Range("D31").Value = "'37.30"
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("D32").Value = "'165.71"
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("D33").Value = "'0.0885"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = "'0.132"
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").Value = "'2.59"
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("D37").Value = "'4.59"
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("D38").Value = "'0.0354"
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = "'3.73"
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("D40").Value = "'2.72"
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("D41").Value = "'1.60"
$ws.Range("E41").Value = '  +6.23%  '
$ws.Range("D42").Value = "'70.46"
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("B43").Value = 'BitcoinSV'
$ws.Range("C43").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D43").Value = "'95.88"
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.229"
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").Value = "'115.20"
$ws.Range("E47").Value = '  +0.99%  '
$ws.Range("D48").Value = "'80.80"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = '1.661.97'
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = "'8.83"
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = "'5.26"
$ws.Range("E51").Value = '  -1.77%  '
